# Refresh cryptos list data (prices / 1h volume change) and restore
# the original row order for WrappedEther/Polygon and Dai/WrappedliquidstakedEther2.0
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.786.53'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.874.59'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("D4").Value = '''1.0000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''0.7294'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '''241.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D7").Value = '''0.9996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.3132'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.64%  '
$ws.Range("D9").Value = '''0.07100'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Value = '''24.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").Value = '''0.08253'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.52%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.906.78'
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.7469'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").Value = '''5.323'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '''92.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '29.790.58'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '''6.031'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '''248.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.21%  '
$ws.Range("D19").Value = '''13.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("D20").Value = '''0.000007814'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''0.9995'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.141.35'
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("D23").Value = '''0.9998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '''7.729'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.40%  '
$ws.Range("D25").Value = '''0.1539'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").Value = '''9.172'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = '''163.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").Value = '''18.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("D29").Value = '''2.028'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").Value = '''1.436'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.07%  '
$ws.Range("D31").Value = '''4.535'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").Value = '''1.528'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").Value = '''4.204'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").Value = '''0.05275'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").Value = '''0.7577'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.05%  '
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("D41").Value = '''0.4485'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.54%  '
$ws.Range("D42").Value = '''6.010'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.97%  '
$ws.Range("D43").Value = '''0.8666'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '''71.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.13%  '
$ws.Range("D45").Value = '1.067.02'
$ws.Range("E45").Value = '  -3.27%  '
$ws.Range("D46").Value = '''104.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.87%  '
$ws.Range("D47").Value = '''0.9998'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").Value = '''1.828'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").Value = '''7.492'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("D50").Value = '''9.518'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").Value = '2.036.28'
$ws.Range("E51").Value = '  +0.79%  '
